$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3665.5356
$ws.Range("J116").Value = 4727
$ws.Range("L116").Value = 4727
$ws.Range("N116").Value = -11611
$ws.Range("H132").Value = 5558102
$ws.Range("I132").Value = 6252474
$ws.Range("J132").Value = 3124.75
$ws.Range("K132").Value = 18757422
$ws.Range("L132").Value = 9374.25
$ws.Range("M132").Value = -18754892
$ws.Range("N132").Value = -14434.25
$ws.Range("H135").Value = 984.7368
$ws.Range("I135").Value = 667.0909
$ws.Range("J135").Value = 3081.2
$ws.Range("K135").Value = 6003.8181
$ws.Range("L135").Value = 27730.8
$ws.Range("M135").Value = -3468.8181
$ws.Range("N135").Value = -32800.8
$ws.Range("H137").Value = 4445.391
$ws.Range("I137").Value = 4197.0527
$ws.Range("K137").Value = 12591.1581
$ws.Range("M137").Value = -10041.1581

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5784.922
$ws.Range("I32").Value = 4869.143
$ws.Range("K32").Value = 4869.143
$ws.Range("M32").Value = -4582.143
$ws.Range("H45").Value = 1526.0244
$ws.Range("I45").Value = 989.9091
$ws.Range("K45").Value = 989.9091
$ws.Range("M45").Value = -612.9091
$ws.Range("H122").Value = 2432.842
$ws.Range("I122").Value = 1629.2858
$ws.Range("K122").Value = 4887.857400000001
$ws.Range("M122").Value = -2437.857400000001
$ws.Range("H132").Value = 21741724
$ws.Range("I132").Value = 31251658
$ws.Range("J132").Value = 4732.357
$ws.Range("K132").Value = 93754974
$ws.Range("L132").Value = 14197.071
$ws.Range("M132").Value = -93752444
$ws.Range("N132").Value = -19257.071

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1869.625
$ws.Range("I80").Value = 2999
$ws.Range("J80").Value = 1493.1666
$ws.Range("K80").Value = 2999
$ws.Range("L80").Value = 1493.1666
$ws.Range("M80").Value = -2001
$ws.Range("N80").Value = -3489.1666
$ws.Range("H83").Value = 1869.625
$ws.Range("I83").Value = 2999
$ws.Range("J83").Value = 1493.1666
$ws.Range("K83").Value = 14995
$ws.Range("L83").Value = 7465.833000000001
$ws.Range("M83").Value = -10003
$ws.Range("N83").Value = -17449.833
$ws.Range("H86").Value = 24721.182
$ws.Range("I86").Value = 1185.375
$ws.Range("J86").Value = 87483.336
$ws.Range("K86").Value = 1185.375
$ws.Range("L86").Value = 87483.336
$ws.Range("M86").Value = -62.375
$ws.Range("N86").Value = -89729.336
$ws.Range("H89").Value = 24721.182
$ws.Range("I89").Value = 1185.375
$ws.Range("J89").Value = 87483.336
$ws.Range("K89").Value = 5926.875
$ws.Range("L89").Value = 437416.68
$ws.Range("M89").Value = -310.875
$ws.Range("N89").Value = -448648.68
$ws.Range("H99").Value = 2323.9167
$ws.Range("I99").Value = 1191
$ws.Range("K99").Value = 1191
$ws.Range("M99").Value = 307
$ws.Range("H107").Value = 2997.5
$ws.Range("I107").Value = 2094.5715
$ws.Range("J107").Value = 5104.3335
$ws.Range("K107").Value = 2094.5715
$ws.Range("L107").Value = 5104.3335
$ws.Range("M107").Value = -174.5715
$ws.Range("N107").Value = -8944.333500000001
$ws.Range("H134").Value = 1934.1714
$ws.Range("I134").Value = 1123.4
$ws.Range("J134").Value = 6798.8
$ws.Range("K134").Value = 3370.2
$ws.Range("L134").Value = 20396.4
$ws.Range("M134").Value = -835.2000000000003
$ws.Range("N134").Value = -25466.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1892.7273
$ws.Range("I16").Value = 780
$ws.Range("J16").Value = 2004
$ws.Range("K16").Value = 780
$ws.Range("L16").Value = 2004
$ws.Range("M16").Value = -493
$ws.Range("N16").Value = -2578
$ws.Range("H31").Value = 2259.2964
$ws.Range("I31").Value = 1739.8368
$ws.Range("J31").Value = 7350
$ws.Range("K31").Value = 1739.8368
$ws.Range("L31").Value = 7350
$ws.Range("M31").Value = -1444.8368
$ws.Range("N31").Value = -7940
$ws.Range("H34").Value = 2259.2964
$ws.Range("I34").Value = 1739.8368
$ws.Range("J34").Value = 7350
$ws.Range("K34").Value = 1739.8368
$ws.Range("L34").Value = 7350
$ws.Range("M34").Value = -1537.8368
$ws.Range("N34").Value = -7754
$ws.Range("H113").Value = 1892.7273
$ws.Range("I113").Value = 780
$ws.Range("J113").Value = 2004
$ws.Range("K113").Value = 780
$ws.Range("L113").Value = 2004
$ws.Range("M113").Value = 1390
$ws.Range("N113").Value = -6344
$ws.Range("H122").Value = 2223.9062
$ws.Range("I122").Value = 1981.1666
$ws.Range("K122").Value = 5943.4998
$ws.Range("M122").Value = -3493.4998
$ws.Range("H132").Value = 2874.6296
$ws.Range("I132").Value = 2172.1428
$ws.Range("J132").Value = 5333.3335
$ws.Range("K132").Value = 6516.428400000001
$ws.Range("L132").Value = 16000.0005
$ws.Range("M132").Value = -3986.428400000001
$ws.Range("N132").Value = -21060.0005
$ws.Range("H134").Value = 1589.1072
$ws.Range("I134").Value = 904.36365
$ws.Range("J134").Value = 4099.8335
$ws.Range("K134").Value = 2713.09095
$ws.Range("L134").Value = 12299.5005
$ws.Range("M134").Value = -178.0909499999998
$ws.Range("N134").Value = -17369.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 13006.923
$ws.Range("I87").Value = 6997.5
$ws.Range("J87").Value = 15677.777
$ws.Range("K87").Value = 20992.5
$ws.Range("L87").Value = 47033.331
$ws.Range("M87").Value = -19744.5
$ws.Range("N87").Value = -49529.331
$ws.Range("H90").Value = 13006.923
$ws.Range("I90").Value = 6997.5
$ws.Range("J90").Value = 15677.777
$ws.Range("K90").Value = 62977.5
$ws.Range("L90").Value = 141099.993
$ws.Range("M90").Value = -56737.5
$ws.Range("N90").Value = -153579.993
$ws.Range("H114").Value = 879.4
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 879.4
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 2638.2
$ws.Range("M114").Value = ""
$ws.Range("N114").Value = -9146.200000000001
$ws.Range("H131").Value = 1209.4656
$ws.Range("J131").Value = 1053.2264
$ws.Range("L131").Value = 3159.6792
$ws.Range("N131").Value = -13239.6792
$ws.Range("H132").Value = 2538.7778
$ws.Range("I132").Value = 1174.8334
$ws.Range("K132").Value = 10573.5006
$ws.Range("M132").Value = -8043.500599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5175.5557
$ws.Range("I70").Value = 5085.7144
$ws.Range("J70").Value = 5490
$ws.Range("K70").Value = 5085.7144
$ws.Range("L70").Value = 5490
$ws.Range("M70").Value = -4815.7144
$ws.Range("N70").Value = -6030
$ws.Range("H73").Value = 5175.5557
$ws.Range("I73").Value = 5085.7144
$ws.Range("J73").Value = 5490
$ws.Range("K73").Value = 5085.7144
$ws.Range("L73").Value = 5490
$ws.Range("M73").Value = -4149.7144
$ws.Range("N73").Value = -7362
$ws.Range("H80").Value = 4011.1538
$ws.Range("I80").Value = 4156.5
$ws.Range("J80").Value = 3526.6667
$ws.Range("K80").Value = 4156.5
$ws.Range("L80").Value = 3526.6667
$ws.Range("M80").Value = -3158.5
$ws.Range("N80").Value = -5522.6667
$ws.Range("H83").Value = 4011.1538
$ws.Range("I83").Value = 4156.5
$ws.Range("J83").Value = 3526.6667
$ws.Range("K83").Value = 20782.5
$ws.Range("L83").Value = 17633.3335
$ws.Range("M83").Value = -15790.5
$ws.Range("N83").Value = -27617.3335
$ws.Range("H126").Value = 4187.375
$ws.Range("I126").Value = 999.5
$ws.Range("J126").Value = 5250
$ws.Range("K126").Value = 2998.5
$ws.Range("L126").Value = 15750
$ws.Range("M126").Value = -528.5
$ws.Range("N126").Value = -20690

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 718.86365
$ws.Range("I16").Value = 738.8095
$ws.Range("J16").Value = 300
$ws.Range("K16").Value = 738.8095
$ws.Range("L16").Value = 300
$ws.Range("M16").Value = -568.8095
$ws.Range("N16").Value = -640
$ws.Range("H22").Value = 35715216
$ws.Range("I22").Value = 71429040
$ws.Range("J22").Value = 1384.9286
$ws.Range("K22").Value = 71429040
$ws.Range("L22").Value = 1384.9286
$ws.Range("M22").Value = -71428745
$ws.Range("N22").Value = -1974.9286
$ws.Range("H27").Value = 35715216
$ws.Range("I27").Value = 71429040
$ws.Range("J27").Value = 1384.9286
$ws.Range("K27").Value = 71429040
$ws.Range("L27").Value = 1384.9286
$ws.Range("M27").Value = -71428933
$ws.Range("N27").Value = -1598.9286
$ws.Range("H46").Value = 2281.818
$ws.Range("I46").Value = 460
$ws.Range("J46").Value = 3800
$ws.Range("K46").Value = 460
$ws.Range("L46").Value = 3800
$ws.Range("M46").Value = -272
$ws.Range("N46").Value = -4176
$ws.Range("H61").Value = 62502876
$ws.Range("I61").Value = 83335250
$ws.Range("J61").Value = 5751.25
$ws.Range("K61").Value = 83335250
$ws.Range("L61").Value = 5751.25
$ws.Range("M61").Value = -83335048
$ws.Range("N61").Value = -6155.25
$ws.Range("H113").Value = 62502876
$ws.Range("I113").Value = 83335250
$ws.Range("J113").Value = 5751.25
$ws.Range("K113").Value = 83335250
$ws.Range("L113").Value = 5751.25
$ws.Range("M113").Value = -83333080
$ws.Range("N113").Value = -10091.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 784.03125
$ws.Range("I136").Value = 513.8929000000001
$ws.Range("J136").Value = 2675
$ws.Range("K136").Value = 1541.6787
$ws.Range("L136").Value = 8025
$ws.Range("M136").Value = 1008.3213
$ws.Range("N136").Value = -13125
